$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 - existing rows 33..62 shift down to 34..63
$ws.Rows(33).Insert()

# Populate the newly inserted row 33 with the new record
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = (Get-Date -Year 2022 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = 100112021
$ws.Cells.Item(33, 7).Value = "Ají"
$ws.Cells.Item(33, 8).Value = "Inferno"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 130
$ws.Cells.Item(33, 11).Value = 22000
$ws.Cells.Item(33, 12).Value = 23000
$ws.Cells.Item(33, 13).Value = 22500
$ws.Cells.Item(33, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 1500
$ws.Cells.Item(33, 17).Value = 15
$ws.Cells.Item(33, 18).Value = "Hortaliza"
